# Auto-generated Excel COM-interop edit script
# Applies scheduled market-data refresh values to the Leve profit tables
# across all class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9604.916999999999
$ws.Range("J9").Value = 114
$ws.Range("L9").Value = 114
$ws.Range("N9").Value = -452
$ws.Range("H106").Value = 2675.1538
$ws.Range("I106").Value = 3961.8333
$ws.Range("K106").Value = 3961.8333
$ws.Range("M106").Value = -3330.8333
$ws.Range("H112").Value = 3166.389
$ws.Range("J112").Value = 3166.389
$ws.Range("L112").Value = 9499.167000000001
$ws.Range("N112").Value = -11715.167
$ws.Range("H129").Value = 2369.1707
$ws.Range("I129").Value = 1024.4286
$ws.Range("K129").Value = 3073.2858
$ws.Range("M129").Value = 1926.7142
$ws.Range("H132").Value = 4776.6523
$ws.Range("I132").Value = 4926.756
$ws.Range("K132").Value = 14780.268
$ws.Range("M132").Value = -12250.268
$ws.Range("H135").Value = 4566.6895
$ws.Range("I135").Value = 4135.9614
$ws.Range("K135").Value = 37223.6526
$ws.Range("M135").Value = -34688.6526
$ws.Range("H137").Value = 20961.61
$ws.Range("J137").Value = 29716.6
$ws.Range("L137").Value = 89149.79999999999
$ws.Range("N137").Value = -94249.79999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 36077.535
$ws.Range("I2").Value = 49757.094
$ws.Range("K2").Value = 49757.094
$ws.Range("M2").Value = -49644.094
$ws.Range("H16").Value = 669
$ws.Range("I16").Value = 669
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 669
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -382
$ws.Range("H45").Value = 103730.7
$ws.Range("I45").Value = 146615.58
$ws.Range("K45").Value = 146615.58
$ws.Range("M45").Value = -146238.58
$ws.Range("H61").Value = 4439.2
$ws.Range("I61").Value = 4439.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4439.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4227.2
$ws.Range("H74").Value = 640003.4399999999
$ws.Range("I74").Value = 717144.8
$ws.Range("J74").Value = 100014
$ws.Range("K74").Value = 717144.8
$ws.Range("L74").Value = 100014
$ws.Range("M74").Value = -716270.8
$ws.Range("N74").Value = -101762
$ws.Range("H77").Value = 640003.4399999999
$ws.Range("I77").Value = 717144.8
$ws.Range("J77").Value = 100014
$ws.Range("K77").Value = 3585724
$ws.Range("L77").Value = 500070
$ws.Range("M77").Value = -3581356
$ws.Range("N77").Value = -508806
$ws.Range("H80").Value = 79864.164
$ws.Range("J80").Value = 79864.164
$ws.Range("L80").Value = 79864.164
$ws.Range("N80").Value = -81860.164
$ws.Range("H83").Value = 79864.164
$ws.Range("J83").Value = 79864.164
$ws.Range("L83").Value = 239592.492
$ws.Range("N83").Value = -249576.492
$ws.Range("H110").Value = 336444.34
$ws.Range("I110").Value = 419154.34
$ws.Range("K110").Value = 419154.34
$ws.Range("M110").Value = -417109.34
$ws.Range("H116").Value = 36077.535
$ws.Range("I116").Value = 49757.094
$ws.Range("K116").Value = 49757.094
$ws.Range("M116").Value = -47463.094
$ws.Range("H132").Value = 213317.52
$ws.Range("J132").Value = 4334.6665
$ws.Range("L132").Value = 13003.9995
$ws.Range("N132").Value = -18063.9995
$ws.Range("H136").Value = 4439.2
$ws.Range("I136").Value = 4439.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13317.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10767.6
$ws.Range("N16").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 36077.535
$ws.Range("I3").Value = 49757.094
$ws.Range("K3").Value = 49757.094
$ws.Range("M3").Value = -49643.094
$ws.Range("H20").Value = 1978.5938
$ws.Range("I20").Value = 1909.3793
$ws.Range("K20").Value = 1909.3793
$ws.Range("M20").Value = -1662.3793
$ws.Range("H134").Value = 32741.473
$ws.Range("I134").Value = 1979.5
$ws.Range("K134").Value = 5938.5
$ws.Range("M134").Value = -3403.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 601391.25
$ws.Range("J31").Value = 26248.277
$ws.Range("L31").Value = 26248.277
$ws.Range("N31").Value = -26838.277
$ws.Range("H34").Value = 601391.25
$ws.Range("J34").Value = 26248.277
$ws.Range("L34").Value = 26248.277
$ws.Range("N34").Value = -26652.277
$ws.Range("H132").Value = 3219
$ws.Range("I132").Value = 1985.4762
$ws.Range("K132").Value = 5956.4286
$ws.Range("M132").Value = -3426.4286
$ws.Range("H134").Value = 209213.4
$ws.Range("I134").Value = 2424.8823
$ws.Range("K134").Value = 7274.646900000001
$ws.Range("M134").Value = -4739.646900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15807984
$ws.Range("I4").Value = 888444.25
$ws.Range("K4").Value = 2665332.75
$ws.Range("M4").Value = -2665220.75
$ws.Range("H88").Value = 7592.3076
$ws.Range("J88").Value = 7592.3076
$ws.Range("L88").Value = 22776.9228
$ws.Range("N88").Value = -23632.9228
$ws.Range("H91").Value = 7592.3076
$ws.Range("J91").Value = 7592.3076
$ws.Range("L91").Value = 22776.9228
$ws.Range("N91").Value = -25740.9228
$ws.Range("H102").Value = 14998.75
$ws.Range("J102").Value = 14998.75
$ws.Range("L102").Value = 44996.25
$ws.Range("N102").Value = -49864.25
$ws.Range("H107").Value = 20369.678
$ws.Range("I107").Value = 666.7143
$ws.Range("J107").Value = 23184.389
$ws.Range("K107").Value = 2000.1429
$ws.Range("L107").Value = 69553.167
$ws.Range("M107").Value = -80.14289999999983
$ws.Range("N107").Value = -73393.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 13532.286
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -212
$ws.Range("H80").Value = 804718.2
$ws.Range("I80").Value = 1003380.9
$ws.Range("J80").Value = 672276.3
$ws.Range("K80").Value = 1003380.9
$ws.Range("L80").Value = 672276.3
$ws.Range("M80").Value = -1002382.9
$ws.Range("N80").Value = -674272.3
$ws.Range("H83").Value = 804718.2
$ws.Range("I83").Value = 1003380.9
$ws.Range("J83").Value = 672276.3
$ws.Range("K83").Value = 5016904.5
$ws.Range("L83").Value = 3361381.5
$ws.Range("M83").Value = -5011912.5
$ws.Range("N83").Value = -3371365.5
$ws.Range("H102").Value = 20408.8
$ws.Range("I102").Value = 28927.846
$ws.Range("K102").Value = 28927.846
$ws.Range("M102").Value = -27305.846
$ws.Range("H122").Value = 508068.88
$ws.Range("I122").Value = 654729.7
$ws.Range("J122").Value = 9422
$ws.Range("K122").Value = 1964189.1
$ws.Range("L122").Value = 28266
$ws.Range("M122").Value = -1961739.1
$ws.Range("N122").Value = -33166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("H46").Value = 4095.52
$ws.Range("I46").Value = 3332.9333
$ws.Range("K46").Value = 3332.9333
$ws.Range("M46").Value = -3144.9333
$ws.Range("H61").Value = 3610.4
$ws.Range("I61").Value = 2747.7896
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 2747.7896
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = -2545.7896
$ws.Range("N61").Value = -20404
$ws.Range("H113").Value = 3610.4
$ws.Range("I113").Value = 2747.7896
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 2747.7896
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = -577.7896000000001
$ws.Range("N113").Value = -24340
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H122").Value = 654871.1
$ws.Range("I122").Value = 5257.4287
$ws.Range("J122").Value = 1109600.8
$ws.Range("K122").Value = 15772.2861
$ws.Range("L122").Value = 3328802.4
$ws.Range("M122").Value = -13322.2861
$ws.Range("N122").Value = -3333702.4
$ws.Range("H123").Value = 84865.5
$ws.Range("J123").Value = 84865.5
$ws.Range("L123").Value = 84865.5
$ws.Range("N123").Value = -94665.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H136").Value = 471768.53
$ws.Range("I136").Value = 532683
$ws.Range("J136").Value = 306429.28
$ws.Range("K136").Value = 1598049
$ws.Range("L136").Value = 919287.8400000001
$ws.Range("M136").Value = -1595499
$ws.Range("N136").Value = -924387.8400000001
$ws.Range("N125").ClearContents()
$ws.Range("N128").ClearContents()
